$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Structural edits first (paragraph deletions / splits), processed from
#    the bottom of the document upward so earlier paragraph indices stay
#    valid as we go.
# ---------------------------------------------------------------------------

# Paragraph 17 (last paragraph): "C" + ".U.14" + ".-Consultar actividades de
# mantenimiento " -- delete entirely.
$p17 = $d.Paragraphs.Item(17)
$p17.Range.Delete()

# Paragraph 16: "C.U.13.-registrar " + [bookmark _GoBack] + "actividad de
# mantenimiento " -- split into two paragraphs: the first keeps the
# (replaced) text, the second keeps only the bookmark.
$p16 = $d.Paragraphs.Item(16)
$r16 = $p16.Range
$r16.Find.Execute("C.U.13.-registrar ", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$r16.InsertParagraphAfter()

$p16a = $d.Paragraphs.Item(16)
$p16a.Range.Find.Execute("C.U.13.-registrar ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "C.U.10.-Registrar Mantenimiento ", 2)

$p16b = $d.Paragraphs.Item(17)
$p16b.Range.Find.Execute("actividad de mantenimiento ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 2)

# Paragraph 13: "C.U.11.-Registrar actividad de mantenimiento " -- delete
# entirely.
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Delete()

# Paragraph 12: "C.U.10.-Registrar licencia" -- delete entirely.
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Delete()

# Paragraph 5: "C.U.3.- " + "Administrar técnicos " (two runs) -- replace the
# first run's text with the final wording and remove the second run's text.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute("C.U.3.- ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "C.U.3.- Administrar Usuario (Técnico Académico) *", 2)
$p5b = $d.Paragraphs.Item(5)
$p5b.Range.Find.Execute("Administrar técnicos ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# 2) Simple single-run text replacements (order independent; each target
#    string is unique in the document).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("C.U. 1.-Administrar hardware ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U. 1.- Administrar Hardware *", 2)

$d.Content.Find.Execute("C.U.2.-Administrar software ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U.2.- Administrar Software *", 2)

$d.Content.Find.Execute("C.U.4.-Administrar usuarios", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U.4.- Administrar Responsable *", 2)

$d.Content.Find.Execute("C.U.5.-Administrar responsable", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U.5.- Generar Datos Estadísticos", 2)

$d.Content.Find.Execute("C.U.6.-Asignar hardware", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CU.6.-  Asignar Hardware", 2)

$d.Content.Find.Execute("C.U.7.-Generar inventario", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U.7.-Generar reporte de fallas e inventario *", 2)

$d.Content.Find.Execute("C.U.8.-Generar datos estadísticos ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U.10.-Administrar licencia", 2)

$d.Content.Find.Execute("C.U.9.-Registrar garantía ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U.10.-Registrar Mantenimiento *", 2)

$d.Content.Find.Execute("C.U.12.-Generar reporte de fallas", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C.U.7.-Generar reporte de fallas", 2)
